$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.155.00"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "'1.935.09"
$ws.Range("E3").Value = "  +2.19%  "
$ws.Range("D4").Value = "'0.9982"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "'326.60"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "'0.9983"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").Value = "'0.4605"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "'0.3898"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "'0.07862"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "'0.9955"
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("D11").Value = "'22.27"
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").Value = "'1.929.42"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("D13").Value = "'5.840"
$ws.Range("E13").Value = "  +2.35%  "
$ws.Range("D14").Value = "'7.100"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "'0.07049"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").Value = "'87.65"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "'0.000009952"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("D20").Value = "'0.9994"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "'29.200.28"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("E22").Value = "  +4.02%  "
$ws.Range("D23").Value = "'11.19"
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("D24").Value = "'2.159.71"
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("D25").Value = "'2.102"
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("D26").Value = "'156.00"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'19.46"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").Value = "'5.899"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "'118.71"
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("D30").Value = "'1.880"
$ws.Range("E30").Value = "  -2.69%  "
$ws.Range("D31").Value = "'0.09326"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "'0.8927"
$ws.Range("E32").Value = "  -2.04%  "
$ws.Range("D33").Value = "'5.212"
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("D34").Value = "'1.326"
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("D35").Value = "'3.147"
$ws.Range("E35").Value = "  -4.12%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "'0.000003524"
$ws.Range("E36").Value = "  +117.97%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.05797"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.169"
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("D39").Value = "'0.02108"
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("D40").Value = "'7.690"
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("D41").Value = "'0.5702"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").Value = "'0.1814"
$ws.Range("E42").Value = "  +2.37%  "
$ws.Range("D43").Value = "'9.751"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").Value = "'11.99"
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("D45").Value = "'2.214"
$ws.Range("E45").Value = "  -1.63%  "
$ws.Range("D46").Value = "'0.5335"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").Value = "'0.06938"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").Value = "'2.595"
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("D49").Value = "'1.848"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "'113.15"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "'0.2985"
$ws.Range("E51").Value = "  +2.52%  "
